# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.653.82'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '1.844.65'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.90'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5272'
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3157'
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06805'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.14'
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7854'
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07777'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '1.842.38'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.38'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.019'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9998'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.94'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007934'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").Value = '26.692.33'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '2.084.40'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.610'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.996'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.353'
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.229'
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.08'
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.684'
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.05'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.06'
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.218'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08705'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.096'
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7308'
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.863'
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.332'
$ws.Range("E38").Value = '  +4.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01739'
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4841'
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9106'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.00'
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.937'
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.719'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4209'
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.106'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1246'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05832'
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.88'
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8967'
$ws.Range("E51").Value = '  +0.99%  '
